$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79, shifting existing rows 79..193 down to 80..194
$ws.Rows("79:79").Insert()

# Populate the newly inserted row 79 with the new weekly record
$ws.Range("A79").Value = 5
$ws.Range("B79").Value = "Macroferia Regional de Talca"
$ws.Range("C79").Value = "Maule"
$ws.Range("D79").Value = 44571
$ws.Range("E79").Value = 7
$ws.Range("F79").Value = 100112008
$ws.Range("G79").Value = "Coliflor"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 2000
$ws.Range("K79").Value = 800
$ws.Range("L79").Value = 800
$ws.Range("M79").Value = 800
$ws.Range("N79").Value = "$/unidad"
$ws.Range("O79").Value = "Región del Maule"
$ws.Range("P79").Value = 800
$ws.Range("Q79").Value = 1
$ws.Range("R79").Value = "Hortaliza"
